# Fix typo in the shared string used by B11: "Olimpio" -> "Olimpico"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Estadio Olimpico Pascual Guerrero"

# Reflect the cursor/selection position recorded in the saved file
$ws.Range("N8").Select()
